$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "I - " + "tipodocumento" (two runs) -> "I - TIPODOCUMENTO" (one run)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("I - tipodocumento", $false, $false, $false, $false, `
    $false, $true, 1, $false, "I - TIPODOCUMENTO", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: split "Vem ao exame desta Comissão o" + " " into
#   "Vem ao exame" / " " / "desta Comissão "
# and move the trailing "o" to after the spell-check markers, directly
# before "pedidoAprovacao":
#   <proofErr spellStart/> "o" "pedidoAprovacao" <proofErr spellEnd/>
# ---------------------------------------------------------------------------
$found = $d.Content
$found.Find.ClearFormatting()
$found.Find.Execute("Vem ao exame desta Comissão o ", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null

# InsertXML needs a "fresh" Range object (not the live Find-mutated range) to
# replace the target span with new run boundaries.
$target = $d.Range($found.Start, $found.End)

$fragment = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
    '<w:r><w:t>Vem ao exame</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">desta Comissão </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>o</w:t></w:r>' + `
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($fragment)
